$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2020" column (Q) is being added to the table, mirroring the
# formatting of the existing "2019" column (P). Copy P3:P13's formatting
# (borders/fonts/number formats) into Q3:Q13 first, then overwrite the
# values: the year heading in Q3, and "-" (no data yet for 2020) in the
# rest of the rows.
$ws.Range("P3:P13").Copy($ws.Range("Q3:Q13")) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("Q3").Value = 2020
$ws.Range("Q4").Value = "-"
$ws.Range("Q5").Value = "-"
$ws.Range("Q6").Value = "-"
$ws.Range("Q7").Value = "-"
$ws.Range("Q8").Value = "-"
$ws.Range("Q9").Value = "-"
$ws.Range("Q10").Value = "-"
$ws.Range("Q11").Value = "-"
$ws.Range("Q12").Value = "-"
$ws.Range("Q13").Value = "-"

# Match the author's final selection / active cell.
$ws.Range("P17").Select() | Out-Null
